$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = "SSPC8D"
$ws.Range("B61").Value = "Film de fusor"
$ws.Range("C61").Value = "HP LaserJet 5000 5200  M5025 M5035 Pro M435 M701 M706, Canon LBP 840 870 910 1610 1810 1820 3500 3900 3950, Canon Image Class 2200 2210 2220 LP3000 LP3010"
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 180000
$ws.Range("F61").Value = 2
$ws.Range("G61").Value = 0
$ws.Range("H61").Formula = "=(E61-D61)*G61"
$ws.Range("I61").Formula = "=D61*F61"
$ws.Range("J61").Value = 0
